$data = @(
  @("ECs","Hras","Agtr1a","ECs",3,1,6.899291333333333,20.697874,0.1714739209386085,0.1714739209386085,2,0.6666666666666666,0.8192526666666667,2.457758,0.03935738656594369,0.03935738656594368,5.652262822943556,50.870365406492,0.006748765392358879,0.006748765392358878),
  @("ECs","Hras","Agtr1a","FAPs",3,1,6.899291333333333,20.697874,0.1714739209386085,0.1714739209386085,3,1,9.221608999999999,27.664827,0.4430115945178313,0.4430115945178313,63.62256705308866,572.603103477798,0.07596493513323749,0.07596493513323747),
  @("ECs","Hras","Agtr1a","sCs",3,1,6.899291333333333,20.697874,0.1714739209386085,0.1714739209386085,3,1,10.77486666666667,32.3246,0.5176310189162251,0.517631018916225,74.33894421115556,669.0504979004,0.08876022041301213,0.08876022041301211),
  @("FAPs","Hras","Agtr1a","ECs",3,1,9.535088,28.605264,0.2369836040920929,0.2369836040920929,2,0.6666666666666666,0.8192526666666667,2.457758,0.03935738656594369,0.03935738656594368,7.811646270901334,70.304816438112,0.009327055316043054,0.009327055316043054),
  @("FAPs","Hras","Agtr1a","FAPs",3,1,9.535088,28.605264,0.2369836040920929,0.2369836040920929,3,1,9.221608999999999,27.664827,0.4430115945178313,0.4430115945178313,87.928853316592,791.3596798493279,0.1049864843234205,0.1049864843234205),
  @("FAPs","Hras","Agtr1a","sCs",3,1,9.535088,28.605264,0.2369836040920929,0.2369836040920929,3,1,10.77486666666667,32.3246,0.5176310189162251,0.517631018916225,102.7393018549333,924.6537166944,0.1226700644526293,0.1226700644526293),
  @("M1","Hras","Agtr1a","ECs",3,1,7.937245333333332,23.811736,0.197271069302819,0.1972710693028191,2,0.6666666666666666,0.8192526666666667,2.457758,0.03935738656594369,0.03935738656594368,6.502609405320888,58.52348464788799,0.007764073732828117,0.007764073732828117),
  @("M1","Hras","Agtr1a","FAPs",3,1,7.937245333333332,23.811736,0.197271069302819,0.1972710693028191,3,1,9.221608999999999,27.664827,0.4430115945178313,0.4430115945178313,73.19417300107465,658.7475570096718,0.08739337096407948,0.08739337096407948),
  @("M1","Hras","Agtr1a","sCs",3,1,7.937245333333332,23.811736,0.197271069302819,0.1972710693028191,3,1,10.77486666666667,32.3246,0.5176310189162251,0.517631018916225,85.52276016728889,769.7048415055999,0.1021136246059115,0.1021136246059115),
  @("M2","Hras","Agtr1a","ECs",3,1,4.424154666666666,13.272464,0.1099572566050275,0.1099572566050275,2,0.6666666666666666,0.8192526666666667,2.457758,0.03935738656594369,0.03935738656594368,3.624500508412444,32.620504575712,0.004327630253934732,0.004327630253934732),
  @("M2","Hras","Agtr1a","FAPs",3,1,4.424154666666666,13.272464,0.1099572566050275,0.1099572566050275,3,1,9.221608999999999,27.664827,0.4430115945178313,0.4430115945178313,40.79782449152533,367.180420423728,0.04871233957739958,0.04871233957739958),
  @("M2","Hras","Agtr1a","sCs",3,1,4.424154666666666,13.272464,0.1099572566050275,0.1099572566050275,3,1,10.77486666666667,32.3246,0.5176310189162251,0.517631018916225,47.66967664604445,429.0270898144,0.05691728677369321,0.0569172867736932),
  @("Neutro","Hras","Agtr1a","ECs",3,1,6.203087666666668,18.609263,0.1541705825626232,0.1541705825626232,2,0.6666666666666666,0.8192526666666667,2.457758,0.03935738656594369,0.03935738656594368,5.081896112483779,45.73706501235401,0.006067751215013899,0.006067751215013899),
  @("Neutro","Hras","Agtr1a","FAPs",3,1,6.203087666666668,18.609263,0.1541705825626232,0.1541705825626232,3,1,9.221608999999999,27.664827,0.4430115945178313,0.4430115945178313,57.20244905472234,514.8220414925011,0.06829935560881066,0.06829935560881067),
  @("Neutro","Hras","Agtr1a","sCs",3,1,6.203087666666668,18.609263,0.1541705825626232,0.1541705825626232,3,1,10.77486666666667,32.3246,0.5176310189162251,0.517631018916225,66.83744252997779,601.5369827698001,0.07980347573879865,0.07980347573879865),
  @("sCs","Hras","Agtr1a","ECs",3,1,5.236355333333333,15.709066,0.1301435664988289,0.1301435664988289,2,0.6666666666666666,0.8192526666666667,2.457758,0.03935738656594369,0.03935738656594368,4.289898070447555,38.609082634028,0.005122110655765008,0.005122110655765008),
  @("sCs","Hras","Agtr1a","FAPs",3,1,5.236355333333333,15.709066,0.1301435664988289,0.1301435664988289,3,1,9.221608999999999,27.664827,0.4430115945178313,0.4430115945178313,48.28762146906466,434.588593221582,0.05765510891088362,0.05765510891088362),
  @("sCs","Hras","Agtr1a","sCs",3,1,5.236355333333333,15.709066,0.1301435664988289,0.1301435664988289,3,1,10.77486666666667,32.3246,0.5176310189162251,0.517631018916225,56.42103053595556,507.7892748236001,0.06736634693218031,0.06736634693218031)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = 2 + $i
  $rowdata = $data[$i]
  for ($j = 0; $j -lt $rowdata.Length; $j++) {
    $c = 1 + $j
    $ws.Cells.Item($r, $c).Value = $rowdata[$j]
  }
}
